$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

Set-TextValue 'D2' '33.969.17'
Set-TextValue 'E2' '  +0.04%  '
Set-TextValue 'D3' '1.781.41'
Set-TextValue 'E3' '  +0.44%  '
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '224.88'
Set-TextValue 'E5' '  +2.21%  '
Set-TextValue 'E6' '  +1.00%  '
Set-TextValue 'E7' '  -0.08%  '
Set-TextValue 'D8' '32.15'
Set-TextValue 'E8' '  +3.55%  '
Set-TextValue 'D9' '0.289'
Set-TextValue 'E9' '  +1.66%  '
Set-TextValue 'D10' '0.0703'
Set-TextValue 'E10' '  +0.62%  '
Set-TextValue 'E11' '  +1.62%  '
Set-TextValue 'D12' '2.034.46'
Set-TextValue 'E12' '  +0.18%  '
Set-TextValue 'B13' 'Chainlink'
Set-TextValue 'C13' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D13' '10.95'
Set-TextValue 'E13' '  +4.56%  '
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.770.60'
Set-TextValue 'E14' '  -0.35%  '
Set-TextValue 'D15' '0.622'
Set-TextValue 'E15' '  +0.55%  '
Set-TextValue 'D16' '33.952.51'
Set-TextValue 'E16' '  +0.06%  '
Set-TextValue 'D17' '4.15'
Set-TextValue 'E17' '  -0.71%  '
Set-TextValue 'D18' '67.79'
Set-TextValue 'E18' '  +0.35%  '
Set-TextValue 'D19' '243.08'
Set-TextValue 'E19' '  +0.44%  '
Set-TextValue 'E20' '  +2.36%  '
Set-TextValue 'E21' '  -0.03%  '
Set-TextValue 'D22' '10.70'
Set-TextValue 'E22' '  +2.35%  '
Set-TextValue 'D23' '4.10'
Set-TextValue 'E23' '  +1.49%  '
Set-TextValue 'E24' '  -2.27%  '
Set-TextValue 'D25' '160.16'
Set-TextValue 'E25' '  +1.74%  '
Set-TextValue 'D26' '16.31'
Set-TextValue 'E26' '  +0.67%  '
Set-TextValue 'E27' '  +1.95%  '
Set-TextValue 'E28' '  +1.33%  '
Set-TextValue 'E29' '  -0.08%  '
Set-TextValue 'D30' '1.24'
Set-TextValue 'E30' '  +4.04%  '
Set-TextValue 'D31' '0.0511'
Set-TextValue 'E31' '  -0.54%  '
Set-TextValue 'D32' '3.65'
Set-TextValue 'E32' '  -0.04%  '
Set-TextValue 'D33' '3.52'
Set-TextValue 'E33' '  +0.81%  '
Set-TextValue 'D34' '1.81'
Set-TextValue 'E34' '  -0.57%  '
Set-TextValue 'D35' '1.393.78'
Set-TextValue 'E35' '  +0.23%  '
Set-TextValue 'D36' '0.658'
Set-TextValue 'E36' '  +5.75%  '
Set-TextValue 'E37' '  -0.85%  '
Set-TextValue 'E38' '  +1.35%  '
Set-TextValue 'E39' '  +7.38%  '
Set-TextValue 'E40' '  +0.84%  '
Set-TextValue 'D41' '0.913'
Set-TextValue 'E41' '  -1.17%  '
Set-TextValue 'B42' 'Aave'
Set-TextValue 'C42' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '77.97'
Set-TextValue 'E42' '  -0.87%  '
Set-TextValue 'B43' 'MXToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D43' '2.67'
Set-TextValue 'E43' '  -0.60%  '
Set-TextValue 'D44' '13.34'
Set-TextValue 'E44' '  +14.89%  '
Set-TextValue 'D45' '0.0₆0139'
Set-TextValue 'E45' '  +18.46%  '
Set-TextValue 'E46' '  +3.77%  '
Set-TextValue 'E47' '  +1.34%  '
Set-TextValue 'D48' '107.96'
Set-TextValue 'E48' '  +3.88%  '
Set-TextValue 'D49' '5.83'
Set-TextValue 'E49' '  +0.25%  '
Set-TextValue 'D50' '1.935.32'
Set-TextValue 'E50' '  +1.07%  '
Set-TextValue 'E51' '  +0.52%  '
